# "Input Consignes et Infos.xlsx" - commit "version stable avec nouveau InputAnnualData"
#
# All real content changes live on the "Consignes-Input-Client" sheet
# (sheet index 2). The layout of the instructions column is reshuffled:
#   - row 2 (B2) gets a corrected wording for the "trimestres" note
#   - the old row 3 (an empty styled A3 cell + the "trimestres RG (...)"
#     note in C3) is removed entirely, which shifts every row below it
#     up by one
#   - the trailing "ATTENTION : si rachat ..." note (old row 11, which
#     becomes row 10 after the shift above) is removed entirely
#   - the selection cursor moves to B4
#
# Sheet "info_régimes" and sheet "Consignes-Input-Inst" keep the exact
# same visible text; only their underlying shared-string ids shift as a
# side effect of the shared string table being rebuilt, which Excel
# handles on its own.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Consignes-Input-Client")

# Fix the wording of the "trimestres" note in B2.
$ws.Range("B2").Value = "trimestres : il s'agit des trimestres tous régimes (pour calcul du taux de retraite et détermination carrière longue)"

# Remove the old row 3 (empty A3 + "trimestres RG (...)" in C3); this
# shifts rows 4-11 up by one, so the old row 4 ("salaires ...") becomes
# the new row 3, old row 5 ("ATTENTION : il faut bien compléter ...")
# becomes row 4, old rows 7-8 (InputCumulDroitsParRegime block) become
# rows 6-7, and old row 10 (InputDateRetraite block) becomes row 9.
$ws.Rows("3").Delete()

# Remove the trailing "ATTENTION : si rachat ..." note, which is now on
# row 10 after the shift above.
$ws.Rows("10").Delete()

# Move the active selection to B4, matching the new layout.
[void]$ws.Range("B4").Select()
